$d = $word.ActiveDocument

# Remove the inline picture (progress-bar image) first.
while ($d.InlineShapes.Count -gt 0) {
    $d.InlineShapes.Item(1).Delete()
}

# The template's body is being reset to a single, empty, centered
# paragraph: strip the three heading paragraphs, the Lorem ipsum
# paragraph, and the blank spacer paragraph that preceded the
# (now-empty) image paragraph.
while ($d.Paragraphs.Count -gt 1) {
    $p1 = $d.Paragraphs.Item(1)
    $p1.Range.Delete()
}
